$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3089511.2
$ws.Range("I76").Value = 3587174.5
$ws.Range("K76").Value = 3587174.5
$ws.Range("M76").Value = -3586859.5
$ws.Range("H79").Value = 3089511.2
$ws.Range("I79").Value = 3587174.5
$ws.Range("K79").Value = 3587174.5
$ws.Range("M79").Value = -3586082.5
$ws.Range("H116").Value = 6593974.5
$ws.Range("J116").Value = 2266.5715
$ws.Range("L116").Value = 2266.5715
$ws.Range("N116").Value = -9150.5715
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1983.6078
$ws.Range("I61").Value = 1514.6046
$ws.Range("J61").Value = 4504.5
$ws.Range("K61").Value = 1514.6046
$ws.Range("L61").Value = 4504.5
$ws.Range("M61").Value = -1302.6046
$ws.Range("N61").Value = -4928.5
$ws.Range("H74").Value = 4561.881
$ws.Range("I74").Value = 1414.697
$ws.Range("J74").Value = 16101.556
$ws.Range("K74").Value = 1414.697
$ws.Range("L74").Value = 16101.556
$ws.Range("M74").Value = -540.6969999999999
$ws.Range("N74").Value = -17849.556
$ws.Range("H77").Value = 4561.881
$ws.Range("I77").Value = 1414.697
$ws.Range("J77").Value = 16101.556
$ws.Range("K77").Value = 7073.485
$ws.Range("L77").Value = 80507.78
$ws.Range("M77").Value = -2705.485
$ws.Range("N77").Value = -89243.78
$ws.Range("H133").Value = 51279.6
$ws.Range("J133").Value = 51279.6
$ws.Range("L133").Value = 51279.6
$ws.Range("N133").Value = -56339.6
$ws.Range("H136").Value = 1983.6078
$ws.Range("I136").Value = 1514.6046
$ws.Range("J136").Value = 4504.5
$ws.Range("K136").Value = 4543.8138
$ws.Range("L136").Value = 13513.5
$ws.Range("M136").Value = -1993.8138
$ws.Range("N136").Value = -18613.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 25974.45
$ws.Range("I16").Value = 34119.2
$ws.Range("J16").Value = 1540.2
$ws.Range("K16").Value = 34119.2
$ws.Range("L16").Value = 1540.2
$ws.Range("M16").Value = -33832.2
$ws.Range("N16").Value = -2114.2
$ws.Range("H31").Value = 1876.1333
$ws.Range("I31").Value = 1086.4412
$ws.Range("K31").Value = 1086.4412
$ws.Range("M31").Value = -791.4412
$ws.Range("H34").Value = 1876.1333
$ws.Range("I34").Value = 1086.4412
$ws.Range("K34").Value = 1086.4412
$ws.Range("M34").Value = -884.4412
$ws.Range("H58").Value = 1119.5555
$ws.Range("I58").Value = 789.1957
$ws.Range("K58").Value = 789.1957
$ws.Range("M58").Value = -586.1957
$ws.Range("H99").Value = 15646167
$ws.Range("I99").Value = 31289444
$ws.Range("J99").Value = 2890
$ws.Range("K99").Value = 31289444
$ws.Range("L99").Value = 2890
$ws.Range("M99").Value = -31287946
$ws.Range("N99").Value = -5886
$ws.Range("H113").Value = 25974.45
$ws.Range("I113").Value = 34119.2
$ws.Range("J113").Value = 1540.2
$ws.Range("K113").Value = 34119.2
$ws.Range("L113").Value = 1540.2
$ws.Range("M113").Value = -31949.2
$ws.Range("N113").Value = -5880.2
$ws.Range("H126").Value = 15646167
$ws.Range("I126").Value = 31289444
$ws.Range("J126").Value = 2890
$ws.Range("K126").Value = 93868332
$ws.Range("L126").Value = 8670
$ws.Range("M126").Value = -93865862
$ws.Range("N126").Value = -13610
$ws.Range("H136").Value = 1119.5555
$ws.Range("I136").Value = 789.1957
$ws.Range("K136").Value = 2367.5871
$ws.Range("M136").Value = 182.4129000000003
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1018.4651
$ws.Range("I5").Value = 627
$ws.Range("J5").Value = 2157.2727
$ws.Range("K5").Value = 1881
$ws.Range("L5").Value = 6471.8181
$ws.Range("M5").Value = -1769
$ws.Range("N5").Value = -6695.8181
$ws.Range("H103").Value = 274938.9
$ws.Range("I103").Value = 1138
$ws.Range("J103").Value = 503106.34
$ws.Range("K103").Value = 3414
$ws.Range("L103").Value = 1509319.02
$ws.Range("M103").Value = -2535
$ws.Range("N103").Value = -1511077.02
$ws.Range("H113").Value = 13889445
$ws.Range("J113").Value = 25000514
$ws.Range("L113").Value = 75001542
$ws.Range("N113").Value = -75005882
$ws.Range("H131").Value = 2143.3713
$ws.Range("I131").Value = 502
$ws.Range("J131").Value = 2416.9333
$ws.Range("K131").Value = 1506
$ws.Range("L131").Value = 7250.7999
$ws.Range("M131").Value = 3534
$ws.Range("N131").Value = -17330.7999
$ws.Range("H135").Value = 1018.4651
$ws.Range("I135").Value = 627
$ws.Range("J135").Value = 2157.2727
$ws.Range("K135").Value = 5643
$ws.Range("L135").Value = 19415.4543
$ws.Range("M135").Value = -3108
$ws.Range("N135").Value = -24485.4543
$ws.Range("H139").Value = 4428
$ws.Range("I139").Value = 4846.6665
$ws.Range("J139").Value = 3800
$ws.Range("K139").Value = 14539.9995
$ws.Range("L139").Value = 11400
$ws.Range("M139").Value = -9399.999500000002
$ws.Range("N139").Value = -21680
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6374.5713
$ws.Range("I70").Value = 6475.684
$ws.Range("J70").Value = 6161.1113
$ws.Range("K70").Value = 6475.684
$ws.Range("L70").Value = 6161.1113
$ws.Range("M70").Value = -6205.684
$ws.Range("N70").Value = -6701.1113
$ws.Range("H73").Value = 6374.5713
$ws.Range("I73").Value = 6475.684
$ws.Range("J73").Value = 6161.1113
$ws.Range("K73").Value = 6475.684
$ws.Range("L73").Value = 6161.1113
$ws.Range("M73").Value = -5539.684
$ws.Range("N73").Value = -8033.1113
$ws.Range("H126").Value = 2261.6667
$ws.Range("I126").Value = 1822.5
$ws.Range("J126").Value = 2481.25
$ws.Range("K126").Value = 5467.5
$ws.Range("L126").Value = 7443.75
$ws.Range("M126").Value = -2997.5
$ws.Range("N126").Value = -12383.75
$ws.Range("H132").Value = 2798.0784
$ws.Range("I132").Value = 2442.8809
$ws.Range("J132").Value = 4455.6665
$ws.Range("K132").Value = 7328.6427
$ws.Range("L132").Value = 13366.9995
$ws.Range("M132").Value = -4798.6427
$ws.Range("N132").Value = -18426.9995
$ws.Range("H138").Value = 65500
$ws.Range("J138").Value = 65500
$ws.Range("L138").Value = 65500
$ws.Range("N138").Value = -75780
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4708.216
$ws.Range("I132").Value = 4498.2812
$ws.Range("J132").Value = 5061.7896
$ws.Range("K132").Value = 13494.8436
$ws.Range("L132").Value = 15185.3688
$ws.Range("M132").Value = -10964.8436
$ws.Range("N132").Value = -20245.3688
$ws.Range("H133").Value = 39437.332
$ws.Range("J133").Value = 39437.332
$ws.Range("L133").Value = 39437.332
$ws.Range("N133").Value = -44497.332
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9093277
$ws.Range("I132").Value = 13891060
$ws.Range("J132").Value = 2743.0527
$ws.Range("K132").Value = 41673180
$ws.Range("L132").Value = 8229.158100000001
$ws.Range("M132").Value = -41670650
$ws.Range("N132").Value = -13289.1581
$ws.Range("H136").Value = 22046.312
$ws.Range("I136").Value = 26245.975
$ws.Range("J136").Value = 3847.7778
$ws.Range("K136").Value = 78737.92499999999
$ws.Range("L136").Value = 11543.3334
$ws.Range("M136").Value = -76187.92499999999
$ws.Range("N136").Value = -16643.3334
